# Horarios actualizados Linea 141 - 156
# Refresh scraped arrival-time tables on all three sheets (LP1912, LP1912-215, 6203-6173):
# new scrape pass at 13:41:54 re-sorted rows by Hora_Llegada and appended newly observed buses.
$wb = $excel.ActiveWorkbook

# ---- LP1912 ----
$ws = $wb.Worksheets("LP1912")
$ws.Range("A2").Value = 'Última actualización: 13:41:54'
$ws.Range("A3").Value = 'Total filas: 191'
$data = New-Object 'object[,]' 191,5
$data[0,0] = '05:57:13'; $data[0,1] = '06:01'; $data[0,2] = '16_SANTA ANA'; $data[0,3] = 4; $data[0,4] = 'LP1912'
$data[1,0] = '05:57:13'; $data[1,1] = '06:09'; $data[1,2] = '10_OLMOS'; $data[1,3] = 12; $data[1,4] = 'LP1912'
$data[2,0] = '05:57:13'; $data[2,1] = '06:16'; $data[2,2] = '215A_EL PATO'; $data[2,3] = 19; $data[2,4] = 'LP1912'
$data[3,0] = '06:17:28'; $data[3,1] = '06:17'; $data[3,2] = '215A_EL PATO'; $data[3,3] = 0; $data[3,4] = 'LP1912'
$data[4,0] = '05:57:13'; $data[4,1] = '06:30'; $data[4,2] = '23_HERNANDEZ'; $data[4,3] = 33; $data[4,4] = 'LP1912'
$data[5,0] = '06:17:28'; $data[5,1] = '06:32'; $data[5,2] = '23_HERNANDEZ'; $data[5,3] = 15; $data[5,4] = 'LP1912'
$data[6,0] = '05:57:13'; $data[6,1] = '06:34'; $data[6,2] = '11_ETCHEVERRY'; $data[6,3] = 37; $data[6,4] = 'LP1912'
$data[7,0] = '06:35:22'; $data[7,1] = '06:35'; $data[7,2] = '11_ETCHEVERRY'; $data[7,3] = 0; $data[7,4] = 'LP1912'
$data[8,0] = '05:57:13'; $data[8,1] = '06:39'; $data[8,2] = '17X38_ROMERO'; $data[8,3] = 42; $data[8,4] = 'LP1912'
$data[9,0] = '05:57:13'; $data[9,1] = '06:41'; $data[9,2] = '16_SANTA ANA'; $data[9,3] = 44; $data[9,4] = 'LP1912'
$data[10,0] = '06:46:50'; $data[10,1] = '06:56'; $data[10,2] = '215A_EL PATO'; $data[10,3] = 10; $data[10,4] = 'LP1912'
$data[11,0] = '05:57:13'; $data[11,1] = '06:57'; $data[11,2] = '215A_EL PATO'; $data[11,3] = 60; $data[11,4] = 'LP1912'
$data[12,0] = '05:57:13'; $data[12,1] = '06:59'; $data[12,2] = '225_GOMEZ'; $data[12,3] = 62; $data[12,4] = 'LP1912'
$data[13,0] = '06:17:28'; $data[13,1] = '07:15'; $data[13,2] = '215C_EL PATO'; $data[13,3] = 58; $data[13,4] = 'LP1912'
$data[14,0] = '05:57:13'; $data[14,1] = '07:16'; $data[14,2] = '215C_EL PATO'; $data[14,3] = 79; $data[14,4] = 'LP1912'
$data[15,0] = '05:57:13'; $data[15,1] = '07:19'; $data[15,2] = '14_ABASTO'; $data[15,3] = 82; $data[15,4] = 'LP1912'
$data[16,0] = '06:46:50'; $data[16,1] = '07:20'; $data[16,2] = '16_SANTA ANA'; $data[16,3] = 34; $data[16,4] = 'LP1912'
$data[17,0] = '05:57:13'; $data[17,1] = '07:21'; $data[17,2] = '23_HERNANDEZ'; $data[17,3] = 84; $data[17,4] = 'LP1912'
$data[18,0] = '06:17:28'; $data[18,1] = '07:21'; $data[18,2] = '16_SANTA ANA'; $data[18,3] = 64; $data[18,4] = 'LP1912'
$data[19,0] = '06:54:04'; $data[19,1] = '07:22'; $data[19,2] = '23_HERNANDEZ'; $data[19,3] = 28; $data[19,4] = 'LP1912'
$data[20,0] = '05:57:13'; $data[20,1] = '07:29'; $data[20,2] = '17X38_ROMERO'; $data[20,3] = 92; $data[20,4] = 'LP1912'
$data[21,0] = '05:57:13'; $data[21,1] = '07:35'; $data[21,2] = '10_OLMOS'; $data[21,3] = 98; $data[21,4] = 'LP1912'
$data[22,0] = '06:17:28'; $data[22,1] = '07:36'; $data[22,2] = '27_EL RETIRO'; $data[22,3] = 79; $data[22,4] = 'LP1912'
$data[23,0] = '05:57:13'; $data[23,1] = '07:37'; $data[23,2] = '27_EL RETIRO'; $data[23,3] = 100; $data[23,4] = 'LP1912'
$data[24,0] = '06:46:50'; $data[24,1] = '07:43'; $data[24,2] = '215A_EL PATO'; $data[24,3] = 57; $data[24,4] = 'LP1912'
$data[25,0] = '06:35:22'; $data[25,1] = '07:44'; $data[25,2] = '215A_EL PATO'; $data[25,3] = 69; $data[25,4] = 'LP1912'
$data[26,0] = '05:57:13'; $data[26,1] = '07:55'; $data[26,2] = '14_ABASTO'; $data[26,3] = 118; $data[26,4] = 'LP1912'
$data[27,0] = '06:17:28'; $data[27,1] = '08:00'; $data[27,2] = '17_ROMERO'; $data[27,3] = 103; $data[27,4] = 'LP1912'
$data[28,0] = '06:46:50'; $data[28,1] = '08:00'; $data[28,2] = '16_SANTA ANA'; $data[28,3] = 74; $data[28,4] = 'LP1912'
$data[29,0] = '06:17:28'; $data[29,1] = '08:01'; $data[29,2] = '16_SANTA ANA'; $data[29,3] = 104; $data[29,4] = 'LP1912'
$data[30,0] = '06:35:22'; $data[30,1] = '08:06'; $data[30,2] = '23_HERNANDEZ'; $data[30,3] = 91; $data[30,4] = 'LP1912'
$data[31,0] = '06:54:04'; $data[31,1] = '08:07'; $data[31,2] = '23_HERNANDEZ'; $data[31,3] = 73; $data[31,4] = 'LP1912'
$data[32,0] = '06:17:28'; $data[32,1] = '08:11'; $data[32,2] = '10_OLMOS'; $data[32,3] = 114; $data[32,4] = 'LP1912'
$data[33,0] = '06:17:28'; $data[33,1] = '08:13'; $data[33,2] = '15X38_ABASTO'; $data[33,3] = 116; $data[33,4] = 'LP1912'
$data[34,0] = '06:35:22'; $data[34,1] = '08:29'; $data[34,2] = '15_ABASTO'; $data[34,3] = 114; $data[34,4] = 'LP1912'
$data[35,0] = '06:35:22'; $data[35,1] = '08:29'; $data[35,2] = '11_ETCHEVERRY'; $data[35,3] = 114; $data[35,4] = 'LP1912'
$data[36,0] = '06:46:50'; $data[36,1] = '08:41'; $data[36,2] = '16_P MOR-SANTA ANA'; $data[36,3] = 115; $data[36,4] = 'LP1912'
$data[37,0] = '06:46:50'; $data[37,1] = '08:43'; $data[37,2] = '215C_EL PATO'; $data[37,3] = 117; $data[37,4] = 'LP1912'
$data[38,0] = '06:54:04'; $data[38,1] = '08:44'; $data[38,2] = '215C_EL PATO'; $data[38,3] = 110; $data[38,4] = 'LP1912'
$data[39,0] = '08:45:36'; $data[39,1] = '08:46'; $data[39,2] = '215C_EL PATO'; $data[39,3] = 1; $data[39,4] = 'LP1912'
$data[40,0] = '07:50:33'; $data[40,1] = '08:49'; $data[40,2] = '23_HERNANDEZ'; $data[40,3] = 59; $data[40,4] = 'LP1912'
$data[41,0] = '07:12:46'; $data[41,1] = '08:51'; $data[41,2] = '23_HERNANDEZ'; $data[41,3] = 99; $data[41,4] = 'LP1912'
$data[42,0] = '08:27:16'; $data[42,1] = '08:52'; $data[42,2] = '23_HERNANDEZ'; $data[42,3] = 25; $data[42,4] = 'LP1912'
$data[43,0] = '07:12:46'; $data[43,1] = '08:53'; $data[43,2] = '215B_EL PATO'; $data[43,3] = 101; $data[43,4] = 'LP1912'
$data[44,0] = '08:10:18'; $data[44,1] = '08:54'; $data[44,2] = '215B_EL PATO'; $data[44,3] = 44; $data[44,4] = 'LP1912'
$data[45,0] = '07:12:46'; $data[45,1] = '08:57'; $data[45,2] = '215A_EL PATO'; $data[45,3] = 105; $data[45,4] = 'LP1912'
$data[46,0] = '07:38:39'; $data[46,1] = '08:58'; $data[46,2] = '215A_EL PATO'; $data[46,3] = 80; $data[46,4] = 'LP1912'
$data[47,0] = '08:10:18'; $data[47,1] = '09:05'; $data[47,2] = '10_OLMOS'; $data[47,3] = 55; $data[47,4] = 'LP1912'
$data[48,0] = '07:38:39'; $data[48,1] = '09:06'; $data[48,2] = '16_SANTA ANA'; $data[48,3] = 88; $data[48,4] = 'LP1912'
$data[49,0] = '07:38:39'; $data[49,1] = '09:16'; $data[49,2] = '27_EL RETIRO'; $data[49,3] = 98; $data[49,4] = 'LP1912'
$data[50,0] = '07:38:39'; $data[50,1] = '09:17'; $data[50,2] = '14_ABASTO'; $data[50,3] = 99; $data[50,4] = 'LP1912'
$data[51,0] = '08:27:16'; $data[51,1] = '09:17'; $data[51,2] = '27_EL RETIRO'; $data[51,3] = 50; $data[51,4] = 'LP1912'
$data[52,0] = '08:10:18'; $data[52,1] = '09:18'; $data[52,2] = '14_ABASTO'; $data[52,3] = 68; $data[52,4] = 'LP1912'
$data[53,0] = '07:38:39'; $data[53,1] = '09:18'; $data[53,2] = '15X38_ABASTO'; $data[53,3] = 100; $data[53,4] = 'LP1912'
$data[54,0] = '07:38:39'; $data[54,1] = '09:29'; $data[54,2] = '10_OLMOS'; $data[54,3] = 111; $data[54,4] = 'LP1912'
$data[55,0] = '08:10:18'; $data[55,1] = '09:31'; $data[55,2] = '16_SANTA ANA'; $data[55,3] = 81; $data[55,4] = 'LP1912'
$data[56,0] = '08:52:50'; $data[56,1] = '09:33'; $data[56,2] = '23_HERNANDEZ'; $data[56,3] = 41; $data[56,4] = 'LP1912'
$data[57,0] = '08:10:18'; $data[57,1] = '09:36'; $data[57,2] = '23_HERNANDEZ'; $data[57,3] = 86; $data[57,4] = 'LP1912'
$data[58,0] = '08:27:16'; $data[58,1] = '09:39'; $data[58,2] = '23_HERNANDEZ'; $data[58,3] = 72; $data[58,4] = 'LP1912'
$data[59,0] = '07:50:33'; $data[59,1] = '09:39'; $data[59,2] = '15_ABASTO'; $data[59,3] = 109; $data[59,4] = 'LP1912'
$data[60,0] = '07:50:33'; $data[60,1] = '09:41'; $data[60,2] = '11_ETCHEVERRY'; $data[60,3] = 111; $data[60,4] = 'LP1912'
$data[61,0] = '08:10:18'; $data[61,1] = '09:42'; $data[61,2] = '11_ETCHEVERRY'; $data[61,3] = 92; $data[61,4] = 'LP1912'
$data[62,0] = '07:50:33'; $data[62,1] = '09:43'; $data[62,2] = '16_P MOR-SANTA ANA'; $data[62,3] = 113; $data[62,4] = 'LP1912'
$data[63,0] = '08:10:18'; $data[63,1] = '09:53'; $data[63,2] = '10_OLMOS'; $data[63,3] = 103; $data[63,4] = 'LP1912'
$data[64,0] = '08:52:50'; $data[64,1] = '09:58'; $data[64,2] = '215C_EL PATO'; $data[64,3] = 66; $data[64,4] = 'LP1912'
$data[65,0] = '08:10:18'; $data[65,1] = '09:59'; $data[65,2] = '215C_EL PATO'; $data[65,3] = 109; $data[65,4] = 'LP1912'
$data[66,0] = '08:37:25'; $data[66,1] = '10:05'; $data[66,2] = '14_ABASTO'; $data[66,3] = 88; $data[66,4] = 'LP1912'
$data[67,0] = '08:10:18'; $data[67,1] = '10:06'; $data[67,2] = '14_ABASTO'; $data[67,3] = 116; $data[67,4] = 'LP1912'
$data[68,0] = '08:27:16'; $data[68,1] = '10:13'; $data[68,2] = '17X38_ROMERO'; $data[68,3] = 106; $data[68,4] = 'LP1912'
$data[69,0] = '09:23:23'; $data[69,1] = '10:21'; $data[69,2] = '23_HERNANDEZ'; $data[69,3] = 58; $data[69,4] = 'LP1912'
$data[70,0] = '09:23:23'; $data[70,1] = '10:25'; $data[70,2] = '16_SANTA ANA'; $data[70,3] = 62; $data[70,4] = 'LP1912'
$data[71,0] = '08:37:25'; $data[71,1] = '10:29'; $data[71,2] = '15_ABASTO'; $data[71,3] = 112; $data[71,4] = 'LP1912'
$data[72,0] = '09:23:23'; $data[72,1] = '10:29'; $data[72,2] = '14_ABASTO'; $data[72,3] = 66; $data[72,4] = 'LP1912'
$data[73,0] = '10:05:51'; $data[73,1] = '10:43'; $data[73,2] = '11X44_ETCHEVERRY'; $data[73,3] = 38; $data[73,4] = 'LP1912'
$data[74,0] = '08:45:36'; $data[74,1] = '10:44'; $data[74,2] = '11X44_ETCHEVERRY'; $data[74,3] = 119; $data[74,4] = 'LP1912'
$data[75,0] = '08:52:50'; $data[75,1] = '10:46'; $data[75,2] = '15_P INDUSTRIAL'; $data[75,3] = 114; $data[75,4] = 'LP1912'
$data[76,0] = '10:05:51'; $data[76,1] = '10:55'; $data[76,2] = '16_SANTA ANA'; $data[76,3] = 50; $data[76,4] = 'LP1912'
$data[77,0] = '10:05:51'; $data[77,1] = '10:56'; $data[77,2] = '27_EL RETIRO'; $data[77,3] = 51; $data[77,4] = 'LP1912'
$data[78,0] = '09:23:23'; $data[78,1] = '10:57'; $data[78,2] = '10_OLMOS'; $data[78,3] = 94; $data[78,4] = 'LP1912'
$data[79,0] = '10:05:51'; $data[79,1] = '10:58'; $data[79,2] = '10_OLMOS'; $data[79,3] = 53; $data[79,4] = 'LP1912'
$data[80,0] = '09:23:23'; $data[80,1] = '10:59'; $data[80,2] = '27_EL RETIRO'; $data[80,3] = 96; $data[80,4] = 'LP1912'
$data[81,0] = '10:50:41'; $data[81,1] = '10:59'; $data[81,2] = '10_OLMOS'; $data[81,3] = 9; $data[81,4] = 'LP1912'
$data[82,0] = '09:23:23'; $data[82,1] = '11:01'; $data[82,2] = '81_EL PELIGRO'; $data[82,3] = 98; $data[82,4] = 'LP1912'
$data[83,0] = '10:05:51'; $data[83,1] = '11:04'; $data[83,2] = '23_HERNANDEZ'; $data[83,3] = 59; $data[83,4] = 'LP1912'
$data[84,0] = '10:37:52'; $data[84,1] = '11:06'; $data[84,2] = '23_HERNANDEZ'; $data[84,3] = 29; $data[84,4] = 'LP1912'
$data[85,0] = '09:23:23'; $data[85,1] = '11:10'; $data[85,2] = '16_P MOR-SANTA ANA'; $data[85,3] = 107; $data[85,4] = 'LP1912'
$data[86,0] = '11:11:33'; $data[86,1] = '11:11'; $data[86,2] = '16_P MOR-SANTA ANA'; $data[86,3] = 0; $data[86,4] = 'LP1912'
$data[87,0] = '09:23:23'; $data[87,1] = '11:14'; $data[87,2] = '14_ABASTO'; $data[87,3] = 111; $data[87,4] = 'LP1912'
$data[88,0] = '09:23:23'; $data[88,1] = '11:15'; $data[88,2] = '15X38_ABASTO'; $data[88,3] = 112; $data[88,4] = 'LP1912'
$data[89,0] = '10:37:52'; $data[89,1] = '11:25'; $data[89,2] = '16_SANTA ANA'; $data[89,3] = 48; $data[89,4] = 'LP1912'
$data[90,0] = '10:05:51'; $data[90,1] = '11:28'; $data[90,2] = '10_OLMOS'; $data[90,3] = 83; $data[90,4] = 'LP1912'
$data[91,0] = '10:50:41'; $data[91,1] = '11:29'; $data[91,2] = '10_OLMOS'; $data[91,3] = 39; $data[91,4] = 'LP1912'
$data[92,0] = '10:05:51'; $data[92,1] = '11:30'; $data[92,2] = '215C_EL PATO'; $data[92,3] = 85; $data[92,4] = 'LP1912'
$data[93,0] = '10:05:51'; $data[93,1] = '11:31'; $data[93,2] = '16_SANTA ANA'; $data[93,3] = 86; $data[93,4] = 'LP1912'
$data[94,0] = '11:11:33'; $data[94,1] = '11:31'; $data[94,2] = '215C_EL PATO'; $data[94,3] = 20; $data[94,4] = 'LP1912'
$data[95,0] = '10:05:51'; $data[95,1] = '11:41'; $data[95,2] = '215B_EL PATO'; $data[95,3] = 96; $data[95,4] = 'LP1912'
$data[96,0] = '11:34:59'; $data[96,1] = '11:44'; $data[96,2] = '15X38_ABASTO'; $data[96,3] = 10; $data[96,4] = 'LP1912'
$data[97,0] = '10:05:51'; $data[97,1] = '11:45'; $data[97,2] = '15X38_ABASTO'; $data[97,3] = 100; $data[97,4] = 'LP1912'
$data[98,0] = '11:47:17'; $data[98,1] = '11:47'; $data[98,2] = '15X38_ABASTO'; $data[98,3] = 0; $data[98,4] = 'LP1912'
$data[99,0] = '11:11:33'; $data[99,1] = '11:51'; $data[99,2] = '23_HERNANDEZ'; $data[99,3] = 40; $data[99,4] = 'LP1912'
$data[100,0] = '11:47:17'; $data[100,1] = '11:52'; $data[100,2] = '23_HERNANDEZ'; $data[100,3] = 5; $data[100,4] = 'LP1912'
$data[101,0] = '11:52:01'; $data[101,1] = '11:52'; $data[101,2] = '15X38_ABASTO'; $data[101,3] = 0; $data[101,4] = 'LP1912'
$data[102,0] = '10:05:51'; $data[102,1] = '11:52'; $data[102,2] = '225_GOMEZ'; $data[102,3] = 107; $data[102,4] = 'LP1912'
$data[103,0] = '10:50:41'; $data[103,1] = '11:53'; $data[103,2] = '225_GOMEZ'; $data[103,3] = 63; $data[103,4] = 'LP1912'
$data[104,0] = '10:37:52'; $data[104,1] = '11:53'; $data[104,2] = '23_HERNANDEZ'; $data[104,3] = 76; $data[104,4] = 'LP1912'
$data[105,0] = '10:50:41'; $data[105,1] = '11:54'; $data[105,2] = '23_HERNANDEZ'; $data[105,3] = 64; $data[105,4] = 'LP1912'
$data[106,0] = '11:52:01'; $data[106,1] = '11:54'; $data[106,2] = '225_GOMEZ'; $data[106,3] = 2; $data[106,4] = 'LP1912'
$data[107,0] = '11:54:18'; $data[107,1] = '11:54'; $data[107,2] = '15X38_ABASTO'; $data[107,3] = 0; $data[107,4] = 'LP1912'
$data[108,0] = '11:34:59'; $data[108,1] = '11:57'; $data[108,2] = '17_ROMERO'; $data[108,3] = 23; $data[108,4] = 'LP1912'
$data[109,0] = '10:05:51'; $data[109,1] = '11:58'; $data[109,2] = '17_ROMERO'; $data[109,3] = 113; $data[109,4] = 'LP1912'
$data[110,0] = '10:37:52'; $data[110,1] = '12:05'; $data[110,2] = '11_ETCHEVERRY'; $data[110,3] = 88; $data[110,4] = 'LP1912'
$data[111,0] = '11:47:17'; $data[111,1] = '12:06'; $data[111,2] = '11_ETCHEVERRY'; $data[111,3] = 19; $data[111,4] = 'LP1912'
$data[112,0] = '11:34:59'; $data[112,1] = '12:09'; $data[112,2] = '15_ABASTO'; $data[112,3] = 35; $data[112,4] = 'LP1912'
$data[113,0] = '11:34:59'; $data[113,1] = '12:09'; $data[113,2] = '16_P MOR-SANTA ANA'; $data[113,3] = 35; $data[113,4] = 'LP1912'
$data[114,0] = '10:37:52'; $data[114,1] = '12:10'; $data[114,2] = '15_ABASTO'; $data[114,3] = 93; $data[114,4] = 'LP1912'
$data[115,0] = '10:37:52'; $data[115,1] = '12:10'; $data[115,2] = '16_P MOR-SANTA ANA'; $data[115,3] = 93; $data[115,4] = 'LP1912'
$data[116,0] = '12:11:52'; $data[116,1] = '12:11'; $data[116,2] = '16_P MOR-SANTA ANA'; $data[116,3] = 0; $data[116,4] = 'LP1912'
$data[117,0] = '12:11:52'; $data[117,1] = '12:13'; $data[117,2] = '15_ABASTO'; $data[117,3] = 2; $data[117,4] = 'LP1912'
$data[118,0] = '10:37:52'; $data[118,1] = '12:16'; $data[118,2] = '10_OLMOS'; $data[118,3] = 99; $data[118,4] = 'LP1912'
$data[119,0] = '11:11:33'; $data[119,1] = '12:17'; $data[119,2] = '10_OLMOS'; $data[119,3] = 66; $data[119,4] = 'LP1912'
$data[120,0] = '10:37:52'; $data[120,1] = '12:21'; $data[120,2] = '215C_EL PATO'; $data[120,3] = 104; $data[120,4] = 'LP1912'
$data[121,0] = '11:11:33'; $data[121,1] = '12:22'; $data[121,2] = '215C_EL PATO'; $data[121,3] = 71; $data[121,4] = 'LP1912'
$data[122,0] = '10:37:52'; $data[122,1] = '12:32'; $data[122,2] = '14_ABASTO'; $data[122,3] = 115; $data[122,4] = 'LP1912'
$data[123,0] = '11:47:17'; $data[123,1] = '12:32'; $data[123,2] = '23_HERNANDEZ'; $data[123,3] = 45; $data[123,4] = 'LP1912'
$data[124,0] = '11:34:59'; $data[124,1] = '12:33'; $data[124,2] = '15_ABASTO'; $data[124,3] = 59; $data[124,4] = 'LP1912'
$data[125,0] = '11:47:17'; $data[125,1] = '12:33'; $data[125,2] = '14_ABASTO'; $data[125,3] = 46; $data[125,4] = 'LP1912'
$data[126,0] = '10:37:52'; $data[126,1] = '12:34'; $data[126,2] = '15_ABASTO'; $data[126,3] = 117; $data[126,4] = 'LP1912'
$data[127,0] = '11:34:59'; $data[127,1] = '12:35'; $data[127,2] = '27_EL RETIRO'; $data[127,3] = 61; $data[127,4] = 'LP1912'
$data[128,0] = '11:11:33'; $data[128,1] = '12:35'; $data[128,2] = '23_HERNANDEZ'; $data[128,3] = 84; $data[128,4] = 'LP1912'
$data[129,0] = '10:50:41'; $data[129,1] = '12:36'; $data[129,2] = '27_EL RETIRO'; $data[129,3] = 106; $data[129,4] = 'LP1912'
$data[130,0] = '11:34:59'; $data[130,1] = '12:36'; $data[130,2] = '23_HERNANDEZ'; $data[130,3] = 62; $data[130,4] = 'LP1912'
$data[131,0] = '11:47:17'; $data[131,1] = '12:37'; $data[131,2] = '27_EL RETIRO'; $data[131,3] = 50; $data[131,4] = 'LP1912'
$data[132,0] = '11:52:01'; $data[132,1] = '12:37'; $data[132,2] = '23_HERNANDEZ'; $data[132,3] = 45; $data[132,4] = 'LP1912'
$data[133,0] = '11:34:59'; $data[133,1] = '12:47'; $data[133,2] = '14_ABASTO'; $data[133,3] = 73; $data[133,4] = 'LP1912'
$data[134,0] = '11:34:59'; $data[134,1] = '12:47'; $data[134,2] = '15X38_ABASTO'; $data[134,3] = 73; $data[134,4] = 'LP1912'
$data[135,0] = '11:34:59'; $data[135,1] = '12:47'; $data[135,2] = '16_SANTA ANA'; $data[135,3] = 73; $data[135,4] = 'LP1912'
$data[136,0] = '10:50:41'; $data[136,1] = '12:48'; $data[136,2] = '16_SANTA ANA'; $data[136,3] = 118; $data[136,4] = 'LP1912'
$data[137,0] = '11:47:17'; $data[137,1] = '12:48'; $data[137,2] = '14_ABASTO'; $data[137,3] = 61; $data[137,4] = 'LP1912'
$data[138,0] = '11:11:33'; $data[138,1] = '12:48'; $data[138,2] = '15X38_ABASTO'; $data[138,3] = 97; $data[138,4] = 'LP1912'
$data[139,0] = '11:11:33'; $data[139,1] = '13:02'; $data[139,2] = '11_ETCHEVERRY'; $data[139,3] = 111; $data[139,4] = 'LP1912'
$data[140,0] = '11:34:59'; $data[140,1] = '13:03'; $data[140,2] = '215C_EL PATO'; $data[140,3] = 89; $data[140,4] = 'LP1912'
$data[141,0] = '11:47:17'; $data[141,1] = '13:03'; $data[141,2] = '11_ETCHEVERRY'; $data[141,3] = 76; $data[141,4] = 'LP1912'
$data[142,0] = '11:47:17'; $data[142,1] = '13:04'; $data[142,2] = '215C_EL PATO'; $data[142,3] = 77; $data[142,4] = 'LP1912'
$data[143,0] = '11:34:59'; $data[143,1] = '13:12'; $data[143,2] = '16_SANTA ANA'; $data[143,3] = 98; $data[143,4] = 'LP1912'
$data[144,0] = '11:47:17'; $data[144,1] = '13:13'; $data[144,2] = '16_SANTA ANA'; $data[144,3] = 86; $data[144,4] = 'LP1912'
$data[145,0] = '11:34:59'; $data[145,1] = '13:16'; $data[145,2] = '10_OLMOS'; $data[145,3] = 102; $data[145,4] = 'LP1912'
$data[146,0] = '11:47:17'; $data[146,1] = '13:17'; $data[146,2] = '10_OLMOS'; $data[146,3] = 90; $data[146,4] = 'LP1912'
$data[147,0] = '12:45:56'; $data[147,1] = '13:18'; $data[147,2] = '15_ABASTO'; $data[147,3] = 33; $data[147,4] = 'LP1912'
$data[148,0] = '12:52:52'; $data[148,1] = '13:19'; $data[148,2] = '15_ABASTO'; $data[148,3] = 27; $data[148,4] = 'LP1912'
$data[149,0] = '12:45:56'; $data[149,1] = '13:21'; $data[149,2] = '23_HERNANDEZ'; $data[149,3] = 36; $data[149,4] = 'LP1912'
$data[150,0] = '11:54:18'; $data[150,1] = '13:22'; $data[150,2] = '23_HERNANDEZ'; $data[150,3] = 88; $data[150,4] = 'LP1912'
$data[151,0] = '11:34:59'; $data[151,1] = '13:24'; $data[151,2] = '16_P MOR-SANTA ANA'; $data[151,3] = 110; $data[151,4] = 'LP1912'
$data[152,0] = '11:47:17'; $data[152,1] = '13:25'; $data[152,2] = '16_P MOR-SANTA ANA'; $data[152,3] = 98; $data[152,4] = 'LP1912'
$data[153,0] = '12:11:52'; $data[153,1] = '13:25'; $data[153,2] = '23_HERNANDEZ'; $data[153,3] = 74; $data[153,4] = 'LP1912'
$data[154,0] = '11:34:59'; $data[154,1] = '13:32'; $data[154,2] = '215A_EL PATO'; $data[154,3] = 118; $data[154,4] = 'LP1912'
$data[155,0] = '12:11:52'; $data[155,1] = '13:32'; $data[155,2] = '14_ABASTO'; $data[155,3] = 81; $data[155,4] = 'LP1912'
$data[156,0] = '11:47:17'; $data[156,1] = '13:33'; $data[156,2] = '215A_EL PATO'; $data[156,3] = 106; $data[156,4] = 'LP1912'
$data[157,0] = '13:41:54'; $data[157,1] = '13:43'; $data[157,2] = '17_ROMERO'; $data[157,3] = 2; $data[157,4] = 'LP1912'
$data[158,0] = '12:11:52'; $data[158,1] = '13:46'; $data[158,2] = '225_GOMEZ'; $data[158,3] = 95; $data[158,4] = 'LP1912'
$data[159,0] = '11:52:01'; $data[159,1] = '13:47'; $data[159,2] = '225_GOMEZ'; $data[159,3] = 115; $data[159,4] = 'LP1912'
$data[160,0] = '12:33:21'; $data[160,1] = '13:54'; $data[160,2] = '15_ABASTO'; $data[160,3] = 81; $data[160,4] = 'LP1912'
$data[161,0] = '12:11:52'; $data[161,1] = '14:01'; $data[161,2] = '10_OLMOS'; $data[161,3] = 110; $data[161,4] = 'LP1912'
$data[162,0] = '12:45:56'; $data[162,1] = '14:01'; $data[162,2] = '23_HERNANDEZ'; $data[162,3] = 76; $data[162,4] = 'LP1912'
$data[163,0] = '12:33:21'; $data[163,1] = '14:02'; $data[163,2] = '10_OLMOS'; $data[163,3] = 89; $data[163,4] = 'LP1912'
$data[164,0] = '13:14:29'; $data[164,1] = '14:02'; $data[164,2] = '16_SANTA ANA'; $data[164,3] = 48; $data[164,4] = 'LP1912'
$data[165,0] = '13:41:54'; $data[165,1] = '14:06'; $data[165,2] = '23_HERNANDEZ'; $data[165,3] = 25; $data[165,4] = 'LP1912'
$data[166,0] = '13:14:29'; $data[166,1] = '14:07'; $data[166,2] = '23_HERNANDEZ'; $data[166,3] = 53; $data[166,4] = 'LP1912'
$data[167,0] = '12:52:52'; $data[167,1] = '14:09'; $data[167,2] = '23_HERNANDEZ'; $data[167,3] = 77; $data[167,4] = 'LP1912'
$data[168,0] = '13:41:54'; $data[168,1] = '14:14'; $data[168,2] = '15_ABASTO'; $data[168,3] = 33; $data[168,4] = 'LP1912'
$data[169,0] = '12:45:56'; $data[169,1] = '14:16'; $data[169,2] = '27_EL RETIRO'; $data[169,3] = 91; $data[169,4] = 'LP1912'
$data[170,0] = '12:33:21'; $data[170,1] = '14:17'; $data[170,2] = '27_EL RETIRO'; $data[170,3] = 104; $data[170,4] = 'LP1912'
$data[171,0] = '12:33:21'; $data[171,1] = '14:17'; $data[171,2] = '11_ETCHEVERRY'; $data[171,3] = 104; $data[171,4] = 'LP1912'
$data[172,0] = '12:45:56'; $data[172,1] = '14:27'; $data[172,2] = '16_SANTA ANA'; $data[172,3] = 102; $data[172,4] = 'LP1912'
$data[173,0] = '13:41:54'; $data[173,1] = '14:28'; $data[173,2] = '16_SANTA ANA'; $data[173,3] = 47; $data[173,4] = 'LP1912'
$data[174,0] = '12:45:56'; $data[174,1] = '14:31'; $data[174,2] = '14X44_ABASTO'; $data[174,3] = 106; $data[174,4] = 'LP1912'
$data[175,0] = '12:33:21'; $data[175,1] = '14:32'; $data[175,2] = '14X44_ABASTO'; $data[175,3] = 119; $data[175,4] = 'LP1912'
$data[176,0] = '12:45:56'; $data[176,1] = '14:33'; $data[176,2] = '215C_EL PATO'; $data[176,3] = 108; $data[176,4] = 'LP1912'
$data[177,0] = '13:14:29'; $data[177,1] = '14:34'; $data[177,2] = '215C_EL PATO'; $data[177,3] = 80; $data[177,4] = 'LP1912'
$data[178,0] = '12:45:56'; $data[178,1] = '14:39'; $data[178,2] = '16_P MOR-SANTA ANA'; $data[178,3] = 114; $data[178,4] = 'LP1912'
$data[179,0] = '12:52:52'; $data[179,1] = '14:47'; $data[179,2] = '215B_EL PATO'; $data[179,3] = 115; $data[179,4] = 'LP1912'
$data[180,0] = '12:52:52'; $data[180,1] = '14:51'; $data[180,2] = '16_SANTA ANA'; $data[180,3] = 119; $data[180,4] = 'LP1912'
$data[181,0] = '13:41:54'; $data[181,1] = '14:51'; $data[181,2] = '23_HERNANDEZ'; $data[181,3] = 70; $data[181,4] = 'LP1912'
$data[182,0] = '13:41:54'; $data[182,1] = '14:53'; $data[182,2] = '215_EL PELIGRO'; $data[182,3] = 72; $data[182,4] = 'LP1912'
$data[183,0] = '13:14:29'; $data[183,1] = '14:54'; $data[183,2] = '215_EL PELIGRO'; $data[183,3] = 100; $data[183,4] = 'LP1912'
$data[184,0] = '13:41:54'; $data[184,1] = '15:01'; $data[184,2] = '10_OLMOS'; $data[184,3] = 80; $data[184,4] = 'LP1912'
$data[185,0] = '13:14:29'; $data[185,1] = '15:02'; $data[185,2] = '10_OLMOS'; $data[185,3] = 108; $data[185,4] = 'LP1912'
$data[186,0] = '13:14:29'; $data[186,1] = '15:12'; $data[186,2] = '14_ABASTO'; $data[186,3] = 118; $data[186,4] = 'LP1912'
$data[187,0] = '13:14:29'; $data[187,1] = '15:13'; $data[187,2] = '17X38_ROMERO'; $data[187,3] = 119; $data[187,4] = 'LP1912'
$data[188,0] = '13:41:54'; $data[188,1] = '15:14'; $data[188,2] = '14_ABASTO'; $data[188,3] = 93; $data[188,4] = 'LP1912'
$data[189,0] = '13:41:54'; $data[189,1] = '15:32'; $data[189,2] = '11_ETCHEVERRY'; $data[189,3] = 111; $data[189,4] = 'LP1912'
$data[190,0] = '13:41:54'; $data[190,1] = '15:33'; $data[190,2] = '215C_EL PATO'; $data[190,3] = 112; $data[190,4] = 'LP1912'
$ws.Range("A6:E196").Value = $data

# ---- LP1912-215 ----
$ws = $wb.Worksheets("LP1912-215")
$ws.Range("A2").Value = 'Última actualización: 13:41:54'
$ws.Range("A3").Value = 'Total filas: 32'
$data = New-Object 'object[,]' 32,5
$data[0,0] = '05:57:13'; $data[0,1] = '06:16'; $data[0,2] = '215A_EL PATO'; $data[0,3] = 19; $data[0,4] = 'LP1912'
$data[1,0] = '06:17:28'; $data[1,1] = '06:17'; $data[1,2] = '215A_EL PATO'; $data[1,3] = 0; $data[1,4] = 'LP1912'
$data[2,0] = '06:46:50'; $data[2,1] = '06:56'; $data[2,2] = '215A_EL PATO'; $data[2,3] = 10; $data[2,4] = 'LP1912'
$data[3,0] = '05:57:13'; $data[3,1] = '06:57'; $data[3,2] = '215A_EL PATO'; $data[3,3] = 60; $data[3,4] = 'LP1912'
$data[4,0] = '06:17:28'; $data[4,1] = '07:15'; $data[4,2] = '215C_EL PATO'; $data[4,3] = 58; $data[4,4] = 'LP1912'
$data[5,0] = '05:57:13'; $data[5,1] = '07:16'; $data[5,2] = '215C_EL PATO'; $data[5,3] = 79; $data[5,4] = 'LP1912'
$data[6,0] = '06:46:50'; $data[6,1] = '07:43'; $data[6,2] = '215A_EL PATO'; $data[6,3] = 57; $data[6,4] = 'LP1912'
$data[7,0] = '06:35:22'; $data[7,1] = '07:44'; $data[7,2] = '215A_EL PATO'; $data[7,3] = 69; $data[7,4] = 'LP1912'
$data[8,0] = '06:46:50'; $data[8,1] = '08:43'; $data[8,2] = '215C_EL PATO'; $data[8,3] = 117; $data[8,4] = 'LP1912'
$data[9,0] = '06:54:04'; $data[9,1] = '08:44'; $data[9,2] = '215C_EL PATO'; $data[9,3] = 110; $data[9,4] = 'LP1912'
$data[10,0] = '08:45:36'; $data[10,1] = '08:46'; $data[10,2] = '215C_EL PATO'; $data[10,3] = 1; $data[10,4] = 'LP1912'
$data[11,0] = '07:12:46'; $data[11,1] = '08:53'; $data[11,2] = '215B_EL PATO'; $data[11,3] = 101; $data[11,4] = 'LP1912'
$data[12,0] = '08:10:18'; $data[12,1] = '08:54'; $data[12,2] = '215B_EL PATO'; $data[12,3] = 44; $data[12,4] = 'LP1912'
$data[13,0] = '07:12:46'; $data[13,1] = '08:57'; $data[13,2] = '215A_EL PATO'; $data[13,3] = 105; $data[13,4] = 'LP1912'
$data[14,0] = '07:38:39'; $data[14,1] = '08:58'; $data[14,2] = '215A_EL PATO'; $data[14,3] = 80; $data[14,4] = 'LP1912'
$data[15,0] = '08:52:50'; $data[15,1] = '09:58'; $data[15,2] = '215C_EL PATO'; $data[15,3] = 66; $data[15,4] = 'LP1912'
$data[16,0] = '08:10:18'; $data[16,1] = '09:59'; $data[16,2] = '215C_EL PATO'; $data[16,3] = 109; $data[16,4] = 'LP1912'
$data[17,0] = '10:05:51'; $data[17,1] = '11:30'; $data[17,2] = '215C_EL PATO'; $data[17,3] = 85; $data[17,4] = 'LP1912'
$data[18,0] = '11:11:33'; $data[18,1] = '11:31'; $data[18,2] = '215C_EL PATO'; $data[18,3] = 20; $data[18,4] = 'LP1912'
$data[19,0] = '10:05:51'; $data[19,1] = '11:41'; $data[19,2] = '215B_EL PATO'; $data[19,3] = 96; $data[19,4] = 'LP1912'
$data[20,0] = '10:37:52'; $data[20,1] = '12:21'; $data[20,2] = '215C_EL PATO'; $data[20,3] = 104; $data[20,4] = 'LP1912'
$data[21,0] = '11:11:33'; $data[21,1] = '12:22'; $data[21,2] = '215C_EL PATO'; $data[21,3] = 71; $data[21,4] = 'LP1912'
$data[22,0] = '11:34:59'; $data[22,1] = '13:03'; $data[22,2] = '215C_EL PATO'; $data[22,3] = 89; $data[22,4] = 'LP1912'
$data[23,0] = '11:47:17'; $data[23,1] = '13:04'; $data[23,2] = '215C_EL PATO'; $data[23,3] = 77; $data[23,4] = 'LP1912'
$data[24,0] = '11:34:59'; $data[24,1] = '13:32'; $data[24,2] = '215A_EL PATO'; $data[24,3] = 118; $data[24,4] = 'LP1912'
$data[25,0] = '11:47:17'; $data[25,1] = '13:33'; $data[25,2] = '215A_EL PATO'; $data[25,3] = 106; $data[25,4] = 'LP1912'
$data[26,0] = '12:45:56'; $data[26,1] = '14:33'; $data[26,2] = '215C_EL PATO'; $data[26,3] = 108; $data[26,4] = 'LP1912'
$data[27,0] = '13:14:29'; $data[27,1] = '14:34'; $data[27,2] = '215C_EL PATO'; $data[27,3] = 80; $data[27,4] = 'LP1912'
$data[28,0] = '12:52:52'; $data[28,1] = '14:47'; $data[28,2] = '215B_EL PATO'; $data[28,3] = 115; $data[28,4] = 'LP1912'
$data[29,0] = '13:41:54'; $data[29,1] = '14:53'; $data[29,2] = '215_EL PELIGRO'; $data[29,3] = 72; $data[29,4] = 'LP1912'
$data[30,0] = '13:14:29'; $data[30,1] = '14:54'; $data[30,2] = '215_EL PELIGRO'; $data[30,3] = 100; $data[30,4] = 'LP1912'
$data[31,0] = '13:41:54'; $data[31,1] = '15:33'; $data[31,2] = '215C_EL PATO'; $data[31,3] = 112; $data[31,4] = 'LP1912'
$ws.Range("A6:E37").Value = $data

# ---- 6203-6173 ----
$ws = $wb.Worksheets("6203-6173")
$ws.Range("A2").Value = 'Última actualización: 13:41:54'
$ws.Range("A3").Value = 'Total filas: 26'
$data = New-Object 'object[,]' 26,5
$data[0,0] = '06:46:50'; $data[0,1] = '07:42'; $data[0,2] = '215A_LA PLATA'; $data[0,3] = 56; $data[0,4] = 'L6173'
$data[1,0] = '05:57:13'; $data[1,1] = '07:43'; $data[1,2] = '215A_LA PLATA'; $data[1,3] = 106; $data[1,4] = 'L6173'
$data[2,0] = '06:46:50'; $data[2,1] = '08:35'; $data[2,2] = '215A_LA PLATA'; $data[2,3] = 109; $data[2,4] = 'L6173'
$data[3,0] = '06:54:04'; $data[3,1] = '08:36'; $data[3,2] = '215A_LA PLATA'; $data[3,3] = 102; $data[3,4] = 'L6173'
$data[4,0] = '08:37:25'; $data[4,1] = '08:37'; $data[4,2] = '215A_LA PLATA'; $data[4,3] = 0; $data[4,4] = 'L6173'
$data[5,0] = '07:12:46'; $data[5,1] = '08:50'; $data[5,2] = '215C_LA PLATA'; $data[5,3] = 98; $data[5,4] = 'L6203'
$data[6,0] = '06:54:04'; $data[6,1] = '08:51'; $data[6,2] = '215C_LA PLATA'; $data[6,3] = 117; $data[6,4] = 'L6203'
$data[7,0] = '08:52:50'; $data[7,1] = '08:52'; $data[7,2] = '215C_LA PLATA'; $data[7,3] = 0; $data[7,4] = 'L6203'
$data[8,0] = '08:52:50'; $data[8,1] = '09:20'; $data[8,2] = '215A_LA PLATA'; $data[8,3] = 28; $data[8,4] = 'L6173'
$data[9,0] = '07:50:33'; $data[9,1] = '09:21'; $data[9,2] = '215A_LA PLATA'; $data[9,3] = 91; $data[9,4] = 'L6173'
$data[10,0] = '08:45:36'; $data[10,1] = '10:12'; $data[10,2] = '215C_LA PLATA'; $data[10,3] = 87; $data[10,4] = 'L6203'
$data[11,0] = '08:27:16'; $data[11,1] = '10:13'; $data[11,2] = '215C_LA PLATA'; $data[11,3] = 106; $data[11,4] = 'L6203'
$data[12,0] = '08:52:50'; $data[12,1] = '10:29'; $data[12,2] = '215B_LP-P MOR-1 Y 57'; $data[12,3] = 97; $data[12,4] = 'L6173'
$data[13,0] = '08:37:25'; $data[13,1] = '10:30'; $data[13,2] = '215B_LP-P MOR-1 Y 57'; $data[13,3] = 113; $data[13,4] = 'L6173'
$data[14,0] = '08:52:50'; $data[14,1] = '10:30'; $data[14,2] = '215A_LA PLATA'; $data[14,3] = 98; $data[14,4] = 'L6173'
$data[15,0] = '08:45:36'; $data[15,1] = '10:31'; $data[15,2] = '215A_LA PLATA'; $data[15,3] = 106; $data[15,4] = 'L6173'
$data[16,0] = '10:05:51'; $data[16,1] = '11:25'; $data[16,2] = '215C_LA PLATA'; $data[16,3] = 80; $data[16,4] = 'L6203'
$data[17,0] = '11:34:59'; $data[17,1] = '13:11'; $data[17,2] = '215C_LA PLATA'; $data[17,3] = 97; $data[17,4] = 'L6203'
$data[18,0] = '11:47:17'; $data[18,1] = '13:12'; $data[18,2] = '215C_LA PLATA'; $data[18,3] = 85; $data[18,4] = 'L6203'
$data[19,0] = '11:34:59'; $data[19,1] = '13:20'; $data[19,2] = '215B_LP-P MOR-1 Y 57'; $data[19,3] = 106; $data[19,4] = 'L6173'
$data[20,0] = '11:47:17'; $data[20,1] = '13:21'; $data[20,2] = '215B_LP-P MOR-1 Y 57'; $data[20,3] = 94; $data[20,4] = 'L6173'
$data[21,0] = '12:11:52'; $data[21,1] = '13:56'; $data[21,2] = '215C_LA PLATA'; $data[21,3] = 105; $data[21,4] = 'L6203'
$data[22,0] = '12:33:21'; $data[22,1] = '13:57'; $data[22,2] = '215C_LA PLATA'; $data[22,3] = 84; $data[22,4] = 'L6203'
$data[23,0] = '13:41:54'; $data[23,1] = '14:28'; $data[23,2] = '215C_LA PLATA'; $data[23,3] = 47; $data[23,4] = 'L6203'
$data[24,0] = '13:14:29'; $data[24,1] = '14:33'; $data[24,2] = '215C_LA PLATA'; $data[24,3] = 79; $data[24,4] = 'L6203'
$data[25,0] = '13:41:54'; $data[25,1] = '15:21'; $data[25,2] = '215A_LA PLATA'; $data[25,3] = 100; $data[25,4] = 'L6173'
$ws.Range("A6:E31").Value = $data

